$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 473
$ws1.Range("F18").Value = 2958
$ws1.Range("F31").Value = 329
$ws1.Range("F32").Value = 1110

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F25").Value = 280
$ws2.Range("F26").Value = 3956

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2476
$ws3.Range("F9").Value = 1339

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2476
$ws4.Range("F7").Value = 1339
$ws4.Range("F11").Value = 473
$ws4.Range("F23").Value = 2958
$ws4.Range("F37").Value = 329
$ws4.Range("F45").Value = 280
$ws4.Range("F46").Value = 1110
